# Apply the "Office Theme" color palette to the presentation, replacing the
# previous "Integral" theme colors (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's ThemeColorScheme indices follow the OOXML <a:clrScheme> child
# order: 1=dk1, 2=lt1, 3=dk2, 4=lt2, 5-10=accent1..accent6, 11=hlink, 12=folHlink.
# ThemeColorScheme.Item(n).RGB is assigned like VBA's RGB(r,g,b) macro, i.e. a
# packed r + g*256 + b*65536 integer, so build that from each target hex code.

function Set-ThemeColor($ColorScheme, $Index, $Hex) {
    $r = [Convert]::ToInt32($Hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($Hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($Hex.Substring(4, 2), 16)
    $packed = $r + ($g * 256) + ($b * 65536)

    $ColorScheme.Item($Index).RGB = $packed
}

$p = $ppt.ActivePresentation

# New "Office Theme" palette, in <a:clrScheme> child order.
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme

for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $index = $i + 1
    $hex = $officeColors[$i]
    Set-ThemeColor $themeColors $index $hex
}
